$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list (price, 1h volume %, and hour columns) for Wed Jan 11
# 2023 16:xx run. Cells keep their existing text storage (NumberFormat
# forced to "@" before the write) so values like "16" or "0.71%" stay
# literal text instead of being reinterpreted as numbers/percentages.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "277.34"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.71%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "16"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.34"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.34%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "16"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.839"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.41%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "16"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06323"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.23%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "16"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.26%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "16"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.417"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.60%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "16"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8899"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.53%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "16"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1531"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.56%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "16"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05279"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5.57%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "16"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07424"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.08%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "16"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02875"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.82%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "16"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08946"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.76%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "16"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001569"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.36%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "16"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006346"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.09%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "16"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006050"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "4.69%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "16"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.470"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.58%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "16"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.302"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.06%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "16"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.13%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "16"

$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "16"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.64%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "16"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.905"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.10%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "16"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "11.44%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "16"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04404"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.18%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "16"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001177"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.37%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "16"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004246"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "10.60%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "16"

$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "16"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001181"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-1.63%"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "16"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0001649"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "-14.85%"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "16"

$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "16"

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "16"

$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "16"

$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "16"

$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "16"

$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "16"

$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "16"

$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "16"

$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "16"

$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "16"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04011"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.57%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "16"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006779"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.72%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "16"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1402"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "19.56%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "16"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001991"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-6.54%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "16"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01170"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.10%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "16"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005360"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.91%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "16"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "9.53%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "16"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01851"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-7.43%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "16"

$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "16"

$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "16"

$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "16"

$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "16"
